{"js": "const replacements = [\n  [\"2025-09-19 Friday\", \"2025-09-20 Saturday\"],\n  [\"87\u00d737=3219\", \"73\u00d743=3139\"],\n  [\"27\u00d732=864\", \"97\u00d752=5044\"],\n  [\"42\u00d784=3528\", \"58\u00d740=2320\"],\n  [\"58\u00d786=4988\", \"31\u00d796=2976\"],\n  [\"61\u00d786=5246\", \"95\u00d727=2565\"],\n  [\"22\u00d742=924\", \"38\u00d781=3078\"],\n  [\"69\u00d754=3726\", \"95\u00d722=2090\"],\n  [\"33\u00d717=561\", \"77\u00d790=6930\"],\n  [\"26\u00d785=2210\", \"94\u00d786=8084\"],\n  [\"87\u00d733=2871\", \"32\u00d736=1152\"],\n  [\"52\u00d723=1196\", \"32\u00d755=1760\"],\n  [\"52\u00d757=2964\", \"44\u00d758=2552\"],\n  [\"56\u00d771=3976\", \"35\u00d739=1365\"],\n  [\"64\u00d797=6208\", \"63\u00d762=3906\"],\n  [\"69\u00d728=1932\", \"23\u00d720=460\"],\n  [\"68\u00d770=4760\", \"84\u00d735=2940\"],\n  [\"67\u00d768=4556\", \"91\u00d751=4641\"],\n  [\"11\u00d733=363\", \"81\u00d798=7938\"],\n  [\"98\u00d745=4410\", \"80\u00d733=2640\"],\n  [\"77\u00d715=1155\", \"46\u00d740=1840\"],\n  [\"14\u00d792=1288\", \"38\u00d767=2546\"],\n  [\"60\u00d797=5820\", \"60\u00d745=2700\"],\n  [\"34\u00d750=1700\", \"96\u00d748=4608\"],\n  [\"32\u00d731=992\", \"74\u00d777=5698\"],\n  [\"21\u00d747=987\", \"52\u00d724=1248\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  for (const [oldText, newText] of replacements) {\n    const results = paragraph.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n    if (results.items.length > 0) {\n      for (const range of results.items) {\n        range.insertText(newText, Word.InsertLocation.replace);\n      }\n      await context.sync();\n      break;\n    }\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-09-19 Friday\", \"2025-09-20 Saturday\"),\n  @(\"87\u00d737=3219\", \"73\u00d743=3139\"),\n  @(\"27\u00d732=864\", \"97\u00d752=5044\"),\n  @(\"42\u00d784=3528\", \"58\u00d740=2320\"),\n  @(\"58\u00d786=4988\", \"31\u00d796=2976\"),\n  @(\"61\u00d786=5246\", \"95\u00d727=2565\"),\n  @(\"22\u00d742=924\", \"38\u00d781=3078\"),\n  @(\"69\u00d754=3726\", \"95\u00d722=2090\"),\n  @(\"33\u00d717=561\", \"77\u00d790=6930\"),\n  @(\"26\u00d785=2210\", \"94\u00d786=8084\"),\n  @(\"87\u00d733=2871\", \"32\u00d736=1152\"),\n  @(\"52\u00d723=1196\", \"32\u00d755=1760\"),\n  @(\"52\u00d757=2964\", \"44\u00d758=2552\"),\n  @(\"56\u00d771=3976\", \"35\u00d739=1365\"),\n  @(\"64\u00d797=6208\", \"63\u00d762=3906\"),\n  @(\"69\u00d728=1932\", \"23\u00d720=460\"),\n  @(\"68\u00d770=4760\", \"84\u00d735=2940\"),\n  @(\"67\u00d768=4556\", \"91\u00d751=4641\"),\n  @(\"11\u00d733=363\", \"81\u00d798=7938\"),\n  @(\"98\u00d745=4410\", \"80\u00d733=2640\"),\n  @(\"77\u00d715=1155\", \"46\u00d740=1840\"),\n  @(\"14\u00d792=1288\", \"38\u00d767=2546\"),\n  @(\"60\u00d797=5820\", \"60\u00d745=2700\"),\n  @(\"34\u00d750=1700\", \"96\u00d748=4608\"),\n  @(\"32\u00d731=992\", \"74\u00d777=5698\"),\n  @(\"21\u00d747=987\", \"52\u00d724=1248\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
